$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.Value = "'41.751.96"
$r.Style = "Normal"
$ws.Range("E2").Value = "  +1.23%  "
$r = $ws.Range("D3")
$r.Value = "'2.268.74"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  -0.05%  "
$r = $ws.Range("D5")
$r.Value = "'304.14"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$r = $ws.Range("D6")
$r.Value = "'91.91"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("E8").Value = "  -0.04%  "
$r = $ws.Range("D9")
$r.Value = "'0.482"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "
$r = $ws.Range("D10")
$r.Value = "'32.29"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +0.31%  "
$r = $ws.Range("D11")
$r.Value = "'53.33"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  -0.02%  "
$r = $ws.Range("D14")
$r.Value = "'6.63"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +0.84%  "
$r = $ws.Range("D15")
$r.Value = "'2.619.72"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +0.76%  "
$r = $ws.Range("D16")
$r.Value = "'14.24"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +1.03%  "
$r = $ws.Range("D17")
$r.Value = "'2.320.83"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +2.63%  "
$r = $ws.Range("D18")
$r.Value = "'0.767"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +2.26%  "
$r = $ws.Range("D19")
$r.Value = "'41.680.11"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +1.29%  "
$r = $ws.Range("D20")
$r.Value = "'12.57"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +6.77%  "
$r = $ws.Range("D21")
$r.Value = "'0.0₃0903"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  +1.50%  "
$r = $ws.Range("D23")
$r.Value = "'66.98"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +0.41%  "
$r = $ws.Range("D24")
$r.Value = "'239.86"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  +0.04%  "
$r = $ws.Range("D27")
$r.Value = "'1.91"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +2.30%  "
$r = $ws.Range("D28")
$r.Value = "'23.89"
$r.Style = "Normal"
$ws.Range("E28").Value = "  -0.05%  "
$r = $ws.Range("D29")
$r.Value = "'9.52"
$r.Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -4.72%  "
$r = $ws.Range("D31")
$r.Value = "'34.78"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +4.78%  "
$r = $ws.Range("D32")
$r.Value = "'160.52"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +0.05%  "
$r = $ws.Range("D33")
$r.Value = "'5.27"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +3.16%  "
$r = $ws.Range("D34")
$r.Value = "'1.00"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("E41").Value = "  +0.69%  "
$r = $ws.Range("D42")
$r.Value = "'3.92"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$r = $ws.Range("D43")
$r.Value = "'2.024.33"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -2.54%  "
$r = $ws.Range("D44")
$r.Value = "'19.22"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -4.83%  "
$r = $ws.Range("D45")
$r.Value = "'10.37"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -0.24%  "
$r = $ws.Range("D46")
$r.Value = "'0.0278"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +0.75%  "
$r = $ws.Range("D47")
$r.Value = "'2.12"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +10.67%  "
$r = $ws.Range("D48")
$r.Value = "'2.88"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$r = $ws.Range("D50")
$r.Value = "'1.15"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$r = $ws.Range("D51")
$r.Value = "'72.30"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +3.33%  "
